$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Realizado (SI/NO)" column for the rows that were marked complete
$ws.Range("C5").Value = "si"
$ws.Range("C7").Value = "si"
$ws.Range("C11").Value = "si"
$ws.Range("C12").Value = "si"
$ws.Range("C13").Value = "si"

# Update the active cell selection
$ws.Range("C8").Select()
